# Populate the "基金受益憑證" (fund) sheet with the full column template
# used by the other asset sheets (name/owner/dealer/quantity/face_value/
# currency/total/property_category/category/date/legislator_name/
# legislator_id/source_file/index), matching row 1 (header) to row 2
# (data) layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("基金受益憑證")

# ---- Row 1: header labels ----
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "owner"
$ws.Range("D1").Value = "dealer"
$ws.Range("E1").Value = "quantity"
$ws.Range("F1").Value = "face_value"
$ws.Range("G1").Value = "currency"
$ws.Range("H1").Value = "total"
$ws.Range("I1").Value = "property_category"
$ws.Range("J1").Value = "category"
$ws.Range("K1").Value = "date"
$ws.Range("L1").Value = "legislator_name"
$ws.Range("M1").Value = "legislator_id"
$ws.Range("N1").Value = "source_file"
$ws.Range("O1").Value = "index"

# Give the newly-added header cells (H1:O1) the same look (bold, border,
# centered) already used by the rest of row 1, without inventing a new
# cell-style entry.
$ws.Range("B1").Copy()
$ws.Range("H1:O1").PasteSpecial(-4122)

# ---- Row 2: data values ----
$ws.Range("A2").Value = 71
$ws.Range("B2").Value = "富達全聚焦"
$ws.Range("C2").Value = "楊際英"
$ws.Range("D2").Value = "台北富邦銀行金華分行"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 38792
$ws.Range("G2").Value = "美金"
$ws.Range("H2").Value = 1171518
$ws.Range("I2").Value = "fund"
$ws.Range("J2").Value = "normal"
$ws.Range("K2").Value = "2011-11-21"
$ws.Range("L2").Value = "蔣乃辛"
$ws.Range("M2").Value = 1722
$ws.Range("N2").Value = "tmp12421"
$ws.Range("O2").Value = 71

# Likewise carry the data-row look onto the newly-added data cells
# (H2:O2) by reusing the existing row-2 style.
$ws.Range("B2").Copy()
$ws.Range("H2:O2").PasteSpecial(-4122)
